$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 4.327115817150455

$ws.Range("B3").Value = 0.1169995834814548
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 4.178645819794754

$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 18.71679738969934
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 24.14949828602258
